# Updating model2 results sheet: refresh R^2 / RMSE / U values and their
# color-scale fills (readme + final models update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function RgbColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# --- Stage a few helper cells (far off to the side) so we can clone the
# exact existing fill formatting (and its palette entry) onto new cells
# via PasteSpecial, instead of creating a brand new fill for colors that
# already exist in the workbook's style table. These helpers are cleared
# again at the end of the script. ---
$ws.Range("H1").Value = 1
$ws.Range("D2").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats; color 00F7FCF5

$ws.Range("H2").Value = 1
$ws.Range("D3").Copy()
$ws.Range("H2").PasteSpecial(-4122) | Out-Null   # color 00B1E0AB

$ws.Range("H3").Value = 1
$ws.Range("E3").Copy()
$ws.Range("H3").PasteSpecial(-4122) | Out-Null   # color 006ABF71

$ws.Range("H5").Value = 1
$ws.Range("E7").Copy()
$ws.Range("H5").PasteSpecial(-4122) | Out-Null   # color 007CC87C

$ws.Range("H6").Value = 1
$ws.Range("E10").Copy()
$ws.Range("H6").PasteSpecial(-4122) | Out-Null   # color 00BCE4B5

$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("C2").Value = 0.7359
$ws.Range("D2").Value = 1.1078
$ws.Range("D2").Interior.Color = RgbColor 231 246 226
$ws.Range("E2").Value = 1.9784
$ws.Range("E2").Interior.Color = RgbColor 232 246 228

# --- Row 3 ---
$ws.Range("C3").Value = 0.8079
$ws.Range("D3").Value = 0.9435
$ws.Range("D3").Interior.Color = RgbColor 180 225 173
$ws.Range("E3").Value = 1.6507
$ws.Range("H2").Copy()
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

# --- Row 4 ---
$ws.Range("C4").Value = 0.8733
$ws.Range("D4").Value = 0.7659
$ws.Range("H3").Copy()
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Value = 1.3273
$ws.Range("E4").Interior.Color = RgbColor 99 188 110

# --- Row 5 (fill/font already correct - dark green w/ white text stays) ---
$ws.Range("C5").Value = 0.9697
$ws.Range("D5").Value = 0.3746
$ws.Range("E5").Value = 0.6389

# --- Row 6 ---
$ws.Range("C6").Value = 0.9143
$ws.Range("D6").Value = 0.6305
$ws.Range("D6").Interior.Color = RgbColor 49 154 80
$ws.Range("E6").Value = 1.1143
$ws.Range("E6").Interior.Color = RgbColor 51 156 82
$ws.Range("E6").Font.Color = 0   # was a dark cell with white text; now lighter, needs black text again

# --- Row 7 ---
$ws.Range("C7").Value = 0.8554
$ws.Range("D7").Value = 0.8077
$ws.Range("H5").Copy()
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Value = 1.4436
$ws.Range("E7").Interior.Color = RgbColor 129 202 129

# --- Row 8 ---
$ws.Range("C8").Value = 0.79
$ws.Range("D8").Value = 0.9639
$ws.Range("H6").Copy()
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = 1.7192
$ws.Range("E8").Interior.Color = RgbColor 192 230 185

# --- Row 9 ---
$ws.Range("C9").Value = 0.7324
$ws.Range("D9").Value = 1.0759
$ws.Range("D9").Interior.Color = RgbColor 222 242 217
$ws.Range("E9").Value = 1.9135
$ws.Range("E9").Interior.Color = RgbColor 224 243 219

# --- Row 10 ---
$ws.Range("C10").Value = 0.657
$ws.Range("D10").Value = 1.2022
$ws.Range("H1").Copy()
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = 2.1327
$ws.Range("H1").Copy()
$ws.Range("E10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Clean up helper cells ---
$ws.Range("H1:H6").Clear()
